$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Handback report generation for 79c843d6-0f05-4992-9646-7c5bdae60136
# (row 8) on both locale sheets: zh-cn and de-de.
# ---------------------------------------------------------------------------

$hyperlinkColor = 15570276   # decimal for RGB(0x64,0x95,0xED) == FF6495ED
$col40 = 39.166666666666664  # ColumnWidth that serializes to width="40"

# ---------------------------- zh-cn sheet ---------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

# New handoff-generation timestamp for this row.
$wsZh.Range("H8").Value = "2016-08-13 00:53:39"

# Latest Target File -- hyperlinked, same as the existing handoff file link.
$wsZh.Hyperlinks.Add($wsZh.Range("I8"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/eeaedeaf5e484e7dcd79d2f8c9ac82cec39eb40c/e2e/79c843d6-0f05-4992-9646-7c5bdae60136.md", "", "", "79c843d6-0f05-4992-9646-7c5bdae60136.md")
$wsZh.Range("I8").Font.Name = "Calibri"
$wsZh.Range("I8").Font.Underline = $true
$wsZh.Range("I8").Font.Color = $hyperlinkColor

# Latest Handback File.
$wsZh.Range("J8").Value = "79c843d6-0f05-4992-9646-7c5bdae60136.036f890c65037f3fd1de0d165eb9cd20e372d300.zh-cn.xlf"

# Latest Handback DateTime -- handback version mismatch error.
$wsZh.Range("K8").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/oltest/blob/2c6aa34f4592aa448a8b682574ec0a4adaf07447/e2e/79c843d6-0f05-4992-9646-7c5bdae60136.md, latest: https://github.com/OpenLocalizationTestOrg/oltest/blob/eeaedeaf5e484e7dcd79d2f8c9ac82cec39eb40c/e2e/79c843d6-0f05-4992-9646-7c5bdae60136.md."

# Error Detail.
$wsZh.Range("P8").Value = "2016-08-13 00:53:11"

# Widen the Error Detail column to fit the new content.
$wsZh.Columns.Item(16).ColumnWidth = $col40

# ---------------------------- de-de sheet ----------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

# Latest Target File -- hyperlinked, same as the existing handoff file link.
$wsDe.Hyperlinks.Add($wsDe.Range("I8"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/eeaedeaf5e484e7dcd79d2f8c9ac82cec39eb40c/e2e/79c843d6-0f05-4992-9646-7c5bdae60136.md", "", "", "79c843d6-0f05-4992-9646-7c5bdae60136.md")
$wsDe.Range("I8").Font.Name = "Calibri"
$wsDe.Range("I8").Font.Underline = $true
$wsDe.Range("I8").Font.Color = $hyperlinkColor

# Latest Handback File.
$wsDe.Range("J8").Value = "79c843d6-0f05-4992-9646-7c5bdae60136.036f890c65037f3fd1de0d165eb9cd20e372d300.de-de.xlf"

# Latest Handback DateTime.
$wsDe.Range("K8").Value = "2016-08-13 00:53:49"

# Error Detail.
$wsDe.Range("P8").Value = "2016-08-13 00:53:11"

# Widen the Error Detail column to fit the new content.
$wsDe.Columns.Item(16).ColumnWidth = $col40
